$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8561180830001831
$ws.Range("B1").Value = 1.423509001731873
$ws.Range("C1").Value = 3.835201025009155
$ws.Range("D1").Value = 2.658177852630615
$ws.Range("E1").Value = 1.60248327255249
